$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)

# New D/E values added to Sheet4 (day-of-month rows), matching the target diff.
# Column D = number of items, Column E = description text.
# Fill order follows the order new shared strings were introduced.

$ws4.Range("D5").Value = 8
$ws4.Range("E5").Value = "Bug Fixes"

$ws4.Range("D8").Value = 3
$ws4.Range("E8").Value = "Bug Fixes"

$ws4.Range("D11").Value = 2
$ws4.Range("E11").Value = "Bug Fixes"

$ws4.Range("D12").Value = 2
$ws4.Range("E12").Value = "Home Page image fix"

$ws4.Range("D17").Value = 2
$ws4.Range("E17").Value = "ETF - date issue fix"

$ws4.Range("D21").Value = 3
$ws4.Range("E21").Value = "Historical report and issue with login in arabic"

$ws4.Range("D25").Value = 4
$ws4.Range("E25").Value = "Factsheet and methodology document fix for indices performance"

$ws4.Range("D26").Value = 2
$ws4.Range("E26").Value = "Basket component feature added"

$ws4.Range("D28").Value = 2
$ws4.Range("E28").Value = "UI Fixes"

$ws4.Range("D29").Value = 8
$ws4.Range("E29").Value = "ETF market watch - new page added"

$ws4.Range("D32").Value = 6
$ws4.Range("E32").Value = "Bug Fixes - Issuer trading info, marketwatch, watchlist"

$ws4.Range("D4").Value = 3
$ws4.Range("E4").Value = "Home page Carousel"

$ws4.Range("D3").Value = 2.5
$ws4.Range("E3").Value = "Forgot password bug fixes"

# Sheet4 becomes the active/selected tab, with G15 as the selected cell.
$ws4.Activate()
$ws4.Range("G15").Select()
